$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update scheduling rows 2-5 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("B2").Value = 46058.20833333334
$wsSchedule.Range("C2").Value = 5
$wsSchedule.Range("D2").Value = 18.9
$wsSchedule.Range("E2").Value = 1361.3674035
$wsSchedule.Range("F2").Value = 72.03002134920635
$wsSchedule.Range("A3").Value = 46058.29166666666
$wsSchedule.Range("C3").Value = 9
$wsSchedule.Range("D3").Value = 34.02
$wsSchedule.Range("E3").Value = 853.5834157500001
$wsSchedule.Range("F3").Value = 25.09063538359789
$wsSchedule.Range("A4").Value = 46058.91666666666
$wsSchedule.Range("B4").Value = 46059.14583333334
$wsSchedule.Range("E4").Value = 1521.22143225
$wsSchedule.Range("F4").Value = 73.17082406204906
$wsSchedule.Range("B5").Value = 46059.64583333334
$wsSchedule.Range("C5").Value = 8.5
$wsSchedule.Range("D5").Value = 32.13
$wsSchedule.Range("E5").Value = 729.14001225
$wsSchedule.Range("F5").Value = 22.69343330999067

# --- Sheet "Detailed": update price/type/status rows 12-97 ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("E12").Value = "OFF"
$wsDetailed.Range("E15").Value = "OFF"
$wsDetailed.Range("B38").Value = 12313.45737
$wsDetailed.Range("B39").Value = 10364.32948
$wsDetailed.Range("B40").Value = 10486.99271
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 1165.90971
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 1129.92235
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 599.6182
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 299.99
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 299.98
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 140.16864
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("E46").Value = "ON"
$wsDetailed.Range("B47").Value = 138.63456
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "ON"
$wsDetailed.Range("B48").Value = 150.88824
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 133.92999
$wsDetailed.Range("B50").Value = 138.42
$wsDetailed.Range("B51").Value = 141.23793
$wsDetailed.Range("B52").Value = 138.42
$wsDetailed.Range("B53").Value = 138.42
$wsDetailed.Range("B54").Value = 138.41996
$wsDetailed.Range("B56").Value = 163.26779
$wsDetailed.Range("B57").Value = 162.21248
$wsDetailed.Range("E57").Value = "OFF"
$wsDetailed.Range("B58").Value = 180.61696
$wsDetailed.Range("E58").Value = "OFF"
$wsDetailed.Range("B59").Value = 197.96695
$wsDetailed.Range("B60").Value = 232.46122
$wsDetailed.Range("B61").Value = 169.81627
$wsDetailed.Range("B62").Value = 267.85792
$wsDetailed.Range("B63").Value = 179.56561
$wsDetailed.Range("B64").Value = 108.89
$wsDetailed.Range("B65").Value = 72.59603
$wsDetailed.Range("B66").Value = 58.78186
$wsDetailed.Range("B67").Value = 56.98
$wsDetailed.Range("B70").Value = 0.51
$wsDetailed.Range("B71").Value = 0.50992
$wsDetailed.Range("B72").Value = 0.50993
$wsDetailed.Range("B73").Value = 0.51
$wsDetailed.Range("B74").Value = 0.51
$wsDetailed.Range("B75").Value = -5.11183
$wsDetailed.Range("B76").Value = 0.51
$wsDetailed.Range("B77").Value = 36.06
$wsDetailed.Range("E79").Value = "ON"
$wsDetailed.Range("B80").Value = 154.2
$wsDetailed.Range("E80").Value = "ON"
$wsDetailed.Range("B82").Value = 919.99
$wsDetailed.Range("B83").Value = 299.75
$wsDetailed.Range("B84").Value = 12224.9073
$wsDetailed.Range("B85").Value = 12289.93225
$wsDetailed.Range("B86").Value = 12301.17287
$wsDetailed.Range("B87").Value = 12328.72822
$wsDetailed.Range("B88").Value = 12927.09746
$wsDetailed.Range("B90").Value = 20222.9
$wsDetailed.Range("B91").Value = 352.7553
$wsDetailed.Range("B92").Value = 299.99
$wsDetailed.Range("B93").Value = 248.84013
$wsDetailed.Range("B94").Value = 207.38775
$wsDetailed.Range("B95").Value = 188.52289
$wsDetailed.Range("B96").Value = 122.25984
$wsDetailed.Range("B97").Value = 138.42
